# Update loading_percent values for the 380 kV case (Case_3_109)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.24730059805024
$ws.Range("C2").Value = 9.051239008829363
$ws.Range("D2").Value = 3.584865416155975
$ws.Range("F2").Value = 16.67961230345439
$ws.Range("G2").Value = 15.79323681549101
$ws.Range("H2").Value = 10.93808153481607
$ws.Range("I2").Value = 14.88288982655263
$ws.Range("N2").Value = 15.41054253696409
$ws.Range("O2").Value = 14.90385827259185

# Row 3
$ws.Range("B3").Value = 10.57380790373071
$ws.Range("C3").Value = 8.68732424265783
$ws.Range("D3").Value = 3.47188645274512
$ws.Range("F3").Value = 16.6726704177993
$ws.Range("G3").Value = 15.751178480952
$ws.Range("H3").Value = 10.97772290876434
$ws.Range("I3").Value = 14.98306941406253
$ws.Range("N3").Value = 15.41033187572917
$ws.Range("O3").Value = 14.95657808676978

# Row 4
$ws.Range("B4").Value = 10.15944325987224
$ws.Range("C4").Value = 8.454693562102944
$ws.Range("D4").Value = 3.399911520074788
$ws.Range("F4").Value = 16.67426726530837
$ws.Range("G4").Value = 15.73364501509494
$ws.Range("H4").Value = 11.00407423107523
$ws.Range("I4").Value = 15.04807860461933
$ws.Range("N4").Value = 15.41179164254423
$ws.Range("O4").Value = 14.99313633327462

# Row 5
$ws.Range("B5").Value = 9.988322807296415
$ws.Range("C5").Value = 8.357682477949993
$ws.Range("D5").Value = 3.369952419810796
$ws.Range("F5").Value = 16.67639140032014
$ws.Range("G5").Value = 15.72858848757138
$ws.Range("H5").Value = 11.0153180435031
$ws.Range("I5").Value = 15.07545068265967
$ws.Range("N5").Value = 15.41278782279457
$ws.Range("O5").Value = 15.00908360703765

# Row 6
$ws.Range("B6").Value = 9.959610190633811
$ws.Range("C6").Value = 8.341443358283204
$ws.Range("D6").Value = 3.364940583563305
$ws.Range("F6").Value = 16.67683307553221
$ws.Range("G6").Value = 15.72787503059339
$ws.Range("H6").Value = 11.01721558707891
$ws.Range("I6").Value = 15.08004898142012
$ws.Range("N6").Value = 15.41297752478997
$ws.Range("O6").Value = 15.01179491667151

# Row 7
$ws.Range("B7").Value = 10.15715554917852
$ws.Range("C7").Value = 8.453394055306973
$ws.Range("D7").Value = 3.399509990065747
$ws.Range("F7").Value = 16.67428994733257
$ws.Range("G7").Value = 15.73356836258142
$ws.Range("H7").Value = 11.00422382302392
$ws.Range("I7").Value = 15.04844418891358
$ws.Range("N7").Value = 15.41180345026411
$ws.Range("O7").Value = 14.99334715936792

# Row 8
$ws.Range("B8").Value = 11.02020395504911
$ws.Range("C8").Value = 8.927722265186668
$ws.Range("D8").Value = 3.546466063578305
$ws.Range("F8").Value = 16.67600313144897
$ws.Range("G8").Value = 15.77701785528259
$ws.Range("H8").Value = 10.95133222256899
$ws.Range("I8").Value = 14.91670597069917
$ws.Range("N8").Value = 15.41014113095308
$ws.Range("O8").Value = 14.92116470225475

# Row 9
$ws.Range("B9").Value = 12.56323696720668
$ws.Range("C9").Value = 9.781334341331947
$ws.Range("D9").Value = 3.812943772417702
$ws.Range("F9").Value = 16.7257850659893
$ws.Range("G9").Value = 15.92764833652888
$ws.Range("H9").Value = 10.86358602710788
$ws.Range("I9").Value = 14.68610052173628
$ws.Range("N9").Value = 15.41941319804421
$ws.Range("O9").Value = 14.81300008004449

# Row 10
$ws.Range("B10").Value = 13.57617667244165
$ws.Range("C10").Value = 10.35770150973811
$ws.Range("D10").Value = 3.994287171559566
$ws.Range("F10").Value = 16.7904738976475
$ws.Range("G10").Value = 16.07745537014242
$ws.Range("H10").Value = 10.80887702460212
$ws.Range("I10").Value = 14.53354439140509
$ws.Range("N10").Value = 15.43375938014443
$ws.Range("O10").Value = 14.7540830456256

# Row 11
$ws.Range("B11").Value = 14.01069014424744
$ws.Range("C11").Value = 10.60820485699964
$ws.Range("D11").Value = 4.073439846326799
$ws.Range("F11").Value = 16.82593627506844
$ws.Range("G11").Value = 16.15386675906674
$ws.Range("H11").Value = 10.78611158307002
$ws.Range("I11").Value = 14.46779898174407
$ws.Range("N11").Value = 15.44189763650564
$ws.Range("O11").Value = 14.73178080206518

# Row 12
$ws.Range("B12").Value = 14.17144867746039
$ws.Range("C12").Value = 10.70133772989664
$ws.Range("D12").Value = 4.102917542720612
$ws.Range("F12").Value = 16.84022444805199
$ws.Range("G12").Value = 16.18396476664281
$ws.Range("H12").Value = 10.77779641052206
$ws.Range("I12").Value = 14.44342795087289
$ws.Range("N12").Value = 15.4452089806005
$ws.Range("O12").Value = 14.7239852670318

# Row 13
$ws.Range("B13").Value = 14.13699468621881
$ws.Range("C13").Value = 10.68135738575153
$ws.Range("D13").Value = 4.096591262778492
$ws.Range("F13").Value = 16.83710916594864
$ws.Range("G13").Value = 16.17743135440377
$ws.Range("H13").Value = 10.7795736323879
$ws.Range("I13").Value = 14.44865332335759
$ws.Range("N13").Value = 15.44448565038713
$ws.Range("O13").Value = 14.72563523095961

# Row 14
$ws.Range("B14").Value = 14.02399165526276
$ws.Range("C14").Value = 10.61590181298404
$ws.Range("D14").Value = 4.075875009435697
$ws.Range("F14").Value = 16.8270946031165
$ws.Range("G14").Value = 16.15631979685303
$ws.Range("H14").Value = 10.78542136135134
$ws.Range("I14").Value = 14.46578343367857
$ws.Range("N14").Value = 15.44216547522794
$ws.Range("O14").Value = 14.73112641968033

# Row 15
$ws.Range("B15").Value = 13.95428142927896
$ws.Range("C15").Value = 10.57558218219637
$ws.Range("D15").Value = 4.063120730604474
$ws.Range("F15").Value = 16.82107202827629
$ws.Range("G15").Value = 16.14353895669618
$ws.Range("H15").Value = 10.78904307695936
$ws.Range("I15").Value = 14.47634453573478
$ws.Range("N15").Value = 15.44077412837013
$ws.Range("O15").Value = 14.73457463915709

# Row 16
$ws.Range("B16").Value = 13.54725102385874
$ws.Range("C16").Value = 10.34109117793176
$ws.Range("D16").Value = 3.989045777816287
$ws.Range("F16").Value = 16.78827709672777
$ws.Range("G16").Value = 16.07262587524058
$ws.Range("H16").Value = 10.81040753943127
$ws.Range("I16").Value = 14.5379145184993
$ws.Range("N16").Value = 15.43325977483831
$ws.Range("O16").Value = 14.75563133745135

# Row 17
$ws.Range("B17").Value = 13.29081686357485
$ws.Range("C17").Value = 10.1942093242255
$ws.Range("D17").Value = 3.942735966001529
$ws.Range("F17").Value = 16.76969862595015
$ws.Range("G17").Value = 16.03122231026193
$ws.Range("H17").Value = 10.82405776737314
$ws.Range("I17").Value = 14.57662135206432
$ws.Range("N17").Value = 15.42906128455169
$ws.Range("O17").Value = 14.76970344570223

# Row 18
$ws.Range("B18").Value = 13.14084925301414
$ws.Range("C18").Value = 10.10862882608225
$ws.Range("D18").Value = 3.915786069350789
$ws.Range("F18").Value = 16.7595814372188
$ws.Range("G18").Value = 16.00818761885961
$ws.Range("H18").Value = 10.83210871851697
$ws.Range("I18").Value = 14.59922834475769
$ws.Range("N18").Value = 15.42679830362652
$ws.Range("O18").Value = 14.7782207033091

# Row 19
$ws.Range("B19").Value = 13.08964814997851
$ws.Range("C19").Value = 10.07946570297922
$ws.Range("D19").Value = 3.906607910626596
$ws.Range("F19").Value = 16.75625383295963
$ws.Range("G19").Value = 16.00052308069538
$ws.Range("H19").Value = 10.8348689166036
$ws.Range("I19").Value = 14.60694174064905
$ws.Range("N19").Value = 15.42605824992598
$ws.Range("O19").Value = 14.78117712107007

# Row 20
$ws.Range("B20").Value = 13.31837074977429
$ws.Range("C20").Value = 10.20995914045576
$ws.Range("D20").Value = 3.947698309287296
$ws.Range("F20").Value = 16.77161753897673
$ws.Range("G20").Value = 16.03554928542177
$ws.Range("H20").Value = 10.82258400555484
$ws.Range("I20").Value = 14.57246535693643
$ws.Range("N20").Value = 15.42949251809676
$ws.Range("O20").Value = 14.76816160893727

# Row 21
$ws.Range("B21").Value = 14.057286044634
$ws.Range("C21").Value = 10.63517491475619
$ws.Range("D21").Value = 4.081973440648522
$ws.Range("F21").Value = 16.83001287576494
$ws.Range("G21").Value = 16.16248943872183
$ws.Range("H21").Value = 10.783695443831
$ws.Range("I21").Value = 14.46073764732469
$ws.Range("N21").Value = 15.44284075505649
$ws.Range("O21").Value = 14.72949586674036

# Row 22
$ws.Range("B22").Value = 14.518164686619
$ws.Range("C22").Value = 10.90299325252001
$ws.Range("D22").Value = 4.166834867392744
$ws.Range("F22").Value = 16.87318152379397
$ws.Range("G22").Value = 16.25221497758075
$ws.Range("H22").Value = 10.7600610390462
$ws.Range("I22").Value = 14.39077920414458
$ws.Range("N22").Value = 15.4529015597162
$ws.Range("O22").Value = 14.70801445675635

# Row 23
$ws.Range("B23").Value = 14.2742024574928
$ws.Range("C23").Value = 10.76098985722705
$ws.Range("D23").Value = 4.121812191725577
$ws.Range("F23").Value = 16.84968693133644
$ws.Range("G23").Value = 16.20371726385156
$ws.Range("H23").Value = 10.77251200327616
$ws.Range("I23").Value = 14.42783713999382
$ws.Range("N23").Value = 15.44741035509898
$ws.Range("O23").Value = 14.71913194137468

# Row 24
$ws.Range("B24").Value = 13.30592154968094
$ws.Range("C24").Value = 10.20284218386042
$ws.Range("D24").Value = 3.945455847416567
$ws.Range("F24").Value = 16.77074824153017
$ws.Range("G24").Value = 16.03359066342212
$ws.Range("H24").Value = 10.82324966019024
$ws.Range("I24").Value = 14.57434317990651
$ws.Range("N24").Value = 15.42929708755401
$ws.Range("O24").Value = 14.76885734344384

# Row 25
$ws.Range("B25").Value = 12.16695921647151
$ws.Range("C25").Value = 9.55906913019013
$ws.Range("D25").Value = 3.743311530967234
$ws.Range("F25").Value = 16.70736121083802
$ws.Range("G25").Value = 15.87995654652285
$ws.Range("H25").Value = 10.88561148933628
$ws.Range("I25").Value = 14.74552047421583
$ws.Range("N25").Value = 15.41557381106028
$ws.Range("O25").Value = 14.83866595189259
